{"js": "// Change \"Design step and Implementation:\" to \"Design steps and Implementation:\"\n// i.e. insert the letter \"s\" right after \"Design step\" (turning \"step\" into \"steps\").\n// The document's trailing \"_GoBack\" bookmark (Word's \"last edit location\" marker)\n// is relocated to sit immediately after the newly inserted \"s\", matching where the\n// edit actually happened, and the now-bookmark-less trailing paragraph is left empty.\n\nconst body = context.document.body;\n\n// Locate the unique heading text and split right after \"Design step\".\nconst hits = body.search(\"Design step\", { matchCase: true, matchWholeWord: false });\nhits.load(\"text\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find 'Design step' in the document body.\");\n}\n\nconst insertionPoint = hits.items[0].getRange(\"End\");\ninsertionPoint.insertText(\"s\", \"Before\");\nawait context.sync();\n\n// Move the \"_GoBack\" bookmark so it sits right after the inserted \"s\" (i.e. right\n// after \"Design steps\"), mirroring where Word leaves it following the edit.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst afterEdit = body.search(\"Design steps\", { matchCase: true });\nafterEdit.load(\"text\");\nawait context.sync();\n\nif (afterEdit.items.length > 0) {\n  const bookmarkSpot = afterEdit.items[0].getRange(\"End\");\n  bookmarkSpot.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Change \"Design step and Implementation:\" to \"Design steps and Implementation:\"\n# i.e. insert the letter \"s\" right after \"Design step\" (turning \"step\" into \"steps\").\n# The document's trailing \"_GoBack\" bookmark (Word's \"last edit location\" marker)\n# is relocated to sit immediately after the newly inserted \"s\", matching where the\n# edit actually happened, and the now-bookmark-less trailing paragraph is left empty.\n\n$d = $word.ActiveDocument\n\n# Find the unique heading text \"Design step\" and collapse to its end.\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"Design step\", $true)\n\nif ($found) {\n    $insertPoint = $d.Range($findRange.End, $findRange.End)\n    $insertPoint.InsertAfter(\"s\")\n}\n\n# Relocate the \"_GoBack\" bookmark to sit right after the inserted \"s\"\n# (i.e. right after \"Design steps\"), mirroring where Word leaves it\n# following the edit.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$findRange2 = $d.Content\n$found2 = $findRange2.Find.Execute(\"Design steps\", $true)\n\nif ($found2) {\n    $bookmarkSpot = $d.Range($findRange2.End, $findRange2.End)\n    $d.Bookmarks.Add(\"_GoBack\", $bookmarkSpot)\n}\n"}
